$wb = $excel.ActiveWorkbook

# Sheet "展览": F4 1323 -> 1324, F5 646 -> 647
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1324
$ws1.Range("F5").Value = 647

# Sheet "全部类型": F4 1323 -> 1324, F6 646 -> 647
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1324
$ws4.Range("F6").Value = 647
